$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.307.98"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "3.522.84"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.15"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.61"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E8").Value = "  +1.83%  "
$ws.Range("E9").Value = "  +8.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.34"
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.437"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "4.115.32"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.51"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("D16").Value = "67.191.95"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "3.516.56"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.17"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "397.51"
$ws.Range("E20").Value = "  +2.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.02"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.59"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.540"
$ws.Range("E23").Value = "  +2.18%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.997"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.31"
$ws.Range("E26").Value = "  +2.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.182"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.31"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.20"
$ws.Range("E32").Value = "  +3.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.40"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("E34").Value = "  +4.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.81"
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.897"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.92"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.77"
$ws.Range("E38").Value = "  +3.57%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.90"
$ws.Range("E39").Value = "  +2.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0749"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.59"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.07"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "2.809.90"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  +3.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.95"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0313"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "342.20"
$ws.Range("E47").Value = "  -4.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.61"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.54"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.853"
$ws.Range("E51").Value = "  +0.38%  "
